# Fix bug with Nonetype in targets field
# - adds a new "т2 конница" unit row's mirrored data into columns N:W of row 6
# - tweaks a few existing numeric values (K4, W5, D6, K6)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Small numeric tweaks
$ws.Range("K4").Value = 5
$ws.Range("W5").Value = 8
$ws.Range("D6").Value = 500
$ws.Range("K6").Value = 1

# Fill in the previously-empty mirrored unit block (N6:W6) for "т2 конница"
$ws.Range("N6").Value = "т2 конница"
$ws.Range("O6").Value = "к"
$ws.Range("P6").Value = 250
$ws.Range("Q6").Value = 250
$ws.Range("R6").Value = 10
$ws.Range("S6").Value = 10
$ws.Range("T6").Value = 0
$ws.Range("U6").Value = 0
$ws.Range("V6").Value = 25
$ws.Range("W6").Value = 2
